$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E12").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("E22").Value = 19
$ws.Range("E23").Value = 16
$ws.Range("E24").Value = 20
$ws.Range("E25").Value = 19
$ws.Range("E26").Value = 22
